$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (2020-04-02 update) as row 45
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A45").Value = 43922
$ws.Range("B45").Value = 7193
$ws.Range("C45").Value = 1677
$ws.Range("D45").Value = 103
$ws.Range("E45").Value = 157
$ws.Range("F45").Value = 5516
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0

# Match the selection state left by the edit
$ws.Range("I45").Select()
